$d = $word.ActiveDocument

# 1) Flujo normal: extend the bullet about the window, adding a search bar and (FA-01)
$d.Content.Find.Execute(
    "un botón “Imprimir” y un botón “Regresar”.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "un botón “Imprimir” y un botón “Regresar” y una barra de búsqueda.  (FA-01)",
    2
) | Out-Null

# 2) Remove the (FA-01) note from the "Imprimir" bullet (now just a trailing space)
$d.Content.Find.Execute(
    "El actor da clic en el botón “Imprimir”. (FA-01)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El actor da clic en el botón “Imprimir”. ",
    2
) | Out-Null

# 3) FA-01 Volver -> FA-01 Clic en barra de "Búsqueda"
$d.Content.Find.Execute(
    "FA-01 Volver",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "FA-01 Clic en barra de “Búsqueda”",
    2
) | Out-Null

# 4) "El actor hace clic en el botón "Regresar"." -> "El actor ingresa información en la búsqueda"
$d.Content.Find.Execute(
    "El actor hace clic en el botón “Regresar”.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El actor ingresa información en la búsqueda",
    2
) | Out-Null

# 5) "El sistema cierra la ventana ReporteView." -> "El sistema filtra la información de PRODUCTOINVENTARIO (nombre) y muestra los resultados obtenidos. "
$d.Content.Find.Execute(
    "El sistema cierra la ventana ReporteView.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El sistema filtra la información de PRODUCTOINVENTARIO (nombre) y muestra los resultados obtenidos. ",
    2
) | Out-Null

# 6) "Fin del caso de uso." (in FA-01 block) -> "Regresa al flujo normal 2. "
$d.Content.Find.Execute(
    "Fin del caso de uso.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Regresa al flujo normal 2. ",
    2
) | Out-Null
